$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '66.912.61'
Set-TextCell 'E2' '  -0.39%  '
Set-TextCell 'D3' '3.115.44'
Set-TextCell 'E3' '  +0.30%  '
Set-TextCell 'E4' '  +0.04%  '
Set-TextCell 'D5' '578.41'
Set-TextCell 'E5' '  -0.72%  '
Set-TextCell 'D6' '172.21'
Set-TextCell 'E6' '  +0.51%  '
Set-TextCell 'D7' '0.999'
Set-TextCell 'E7' '  +0.04%  '
Set-TextCell 'D8' '3.112.83'
Set-TextCell 'E8' '  +0.38%  '
Set-TextCell 'E9' '  -1.17%  '
Set-TextCell 'E10' '  -3.31%  '
Set-TextCell 'E11' '  -1.96%  '
Set-TextCell 'D12' '0.483'
Set-TextCell 'E12' '  +0.00%  '
Set-TextCell 'D13' '0.0000246'
Set-TextCell 'E13' '  -2.50%  '
Set-TextCell 'D14' '37.21'
Set-TextCell 'E14' '  +0.14%  '
Set-TextCell 'E15' '  -1.53%  '
Set-TextCell 'D16' '3.631.71'
Set-TextCell 'E16' '  +0.19%  '
Set-TextCell 'D17' '66.843.01'
Set-TextCell 'E17' '  -0.52%  '
Set-TextCell 'E18' '  -1.56%  '
Set-TextCell 'D19' '3.113.41'
Set-TextCell 'E19' '  +0.10%  '
Set-TextCell 'E20' '  +1.55%  '
Set-TextCell 'D21' '476.80'
Set-TextCell 'E21' '  +0.65%  '
Set-TextCell 'B22' 'Uniswap'
Set-TextCell 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D22' '8.03'
Set-TextCell 'E22' '  +6.10%  '
Set-TextCell 'B23' 'Polygon'
Set-TextCell 'C23' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D23' '0.714'
Set-TextCell 'E23' '  -0.79%  '
Set-TextCell 'D24' '13.55'
Set-TextCell 'E24' '  +4.24%  '
Set-TextCell 'D25' '84.03'
Set-TextCell 'E25' '  -0.01%  '
Set-TextCell 'D26' '2.30'
Set-TextCell 'E26' '  -3.31%  '
Set-TextCell 'B27' 'Dai'
Set-TextCell 'C27' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D27' '1.00'
Set-TextCell 'E27' '  -0.04%  '
Set-TextCell 'B28' 'RenderToken'
Set-TextCell 'C28' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D28' '10.00'
Set-TextCell 'E28' '  -3.24%  '
Set-TextCell 'D29' '2.40'
Set-TextCell 'E29' '  -1.88%  '
Set-TextCell 'D30' '7.89'
Set-TextCell 'E30' '  -3.19%  '
Set-TextCell 'D31' '2.67'
Set-TextCell 'E31' '  -0.78%  '
Set-TextCell 'D32' '28.62'
Set-TextCell 'E32' '  +0.24%  '
Set-TextCell 'E33' '  -0.74%  '
Set-TextCell 'D34' '0.0₃0944'
Set-TextCell 'E34' '  -8.26%  '
Set-TextCell 'E35' '  +0.05%  '
Set-TextCell 'D36' '5.87'
Set-TextCell 'E36' '  -1.39%  '
Set-TextCell 'D37' '0.979'
Set-TextCell 'E37' '  -3.73%  '
Set-TextCell 'D38' '47.21'
Set-TextCell 'E38' '  -0.25%  '
Set-TextCell 'D39' '2.08'
Set-TextCell 'E39' '  -2.19%  '
Set-TextCell 'D40' '50.03'
Set-TextCell 'E40' '  -0.96%  '
Set-TextCell 'D41' '0.310'
Set-TextCell 'E41' '  -3.37%  '
Set-TextCell 'E42' '  -1.75%  '
Set-TextCell 'D43' '8.68'
Set-TextCell 'E43' '  -0.76%  '
Set-TextCell 'D44' '2.813.39'
Set-TextCell 'E44' '  +1.23%  '
Set-TextCell 'B45' 'dogwifhat'
Set-TextCell 'C45' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D45' '2.60'
Set-TextCell 'E45' '  -11.05%  '
Set-TextCell 'B46' 'VeChain'
Set-TextCell 'C46' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D46' '0.0356'
Set-TextCell 'E46' '  -2.59%  '
Set-TextCell 'B47' 'Bittensor'
Set-TextCell 'C47' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D47' '381.11'
Set-TextCell 'E47' '  -4.00%  '
Set-TextCell 'D48' '136.06'
Set-TextCell 'E48' '  +0.68%  '
Set-TextCell 'E49' '  +0.03%  '
Set-TextCell 'D50' '24.83'
Set-TextCell 'E50' '  -0.23%  '
Set-TextCell 'E51' '  -2.52%  '
